# Apply updated cryptocurrency price/volume data per the Dec 15 2023 GitHub Actions refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.231.98'
$ws.Range('E2').Value = '  -1.79%  '
$ws.Range('D3').Value = '2.247.74'
$ws.Range('E3').Value = '  -2.27%  '
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').Value = '''247.49'
$ws.Range('E5').Value = '  -2.10%  '
$ws.Range('D6').Value = '''0.634'
$ws.Range('E6').Value = '  -1.48%  '
$ws.Range('D7').Value = '''77.03'
$ws.Range('E7').Value = '  +3.57%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = '''0.627'
$ws.Range('E9').Value = '  -2.85%  '
$ws.Range('D10').Value = '''42.34'
$ws.Range('E10').Value = '  +7.00%  '
$ws.Range('E11').Value = '  -3.35%  '
$ws.Range('D12').Value = '''7.12'
$ws.Range('E12').Value = '  -5.17%  '
$ws.Range('E13').Value = '  -3.25%  '
$ws.Range('D14').Value = '2.582.93'
$ws.Range('E14').Value = '  -2.33%  '
$ws.Range('D15').Value = '''14.70'
$ws.Range('E15').Value = '  -4.51%  '
$ws.Range('D16').Value = '''0.858'
$ws.Range('E16').Value = '  -2.06%  '
$ws.Range('D17').Value = '2.248.03'
$ws.Range('E17').Value = '  -2.39%  '
$ws.Range('D18').Value = '42.088.16'
$ws.Range('E18').Value = '  -1.97%  '
$ws.Range('E19').Value = '  -2.90%  '
$ws.Range('E20').Value = '  -1.20%  '
$ws.Range('E21').Value = '  -3.45%  '
$ws.Range('D22').Value = '''2.28'
$ws.Range('E22').Value = '  +0.21%  '
$ws.Range('D23').Value = '''232.06'
$ws.Range('E23').Value = '  -2.85%  '
$ws.Range('E24').Value = '  -0.03%  '
$ws.Range('D25').Value = '''11.34'
$ws.Range('E25').Value = '  -2.62%  '
$ws.Range('E26').Value = '  -7.40%  '
$ws.Range('E27').Value = '  -5.13%  '
$ws.Range('D28').Value = '''7.44'
$ws.Range('E28').Value = '  +17.92%  '
$ws.Range('D29').Value = '''2.15'
$ws.Range('E29').Value = '  +0.22%  '
$ws.Range('D30').Value = '''169.78'
$ws.Range('E30').Value = '  +1.25%  '
$ws.Range('D31').Value = '''20.59'
$ws.Range('E31').Value = '  -2.71%  '
$ws.Range('D32').Value = '''0.0836'
$ws.Range('E32').Value = '  -0.61%  '
$ws.Range('D33').Value = '''32.62'
$ws.Range('E33').Value = '  +5.31%  '
$ws.Range('D34').Value = '''0.120'
$ws.Range('E34').Value = '  -5.15%  '
$ws.Range('D35').Value = '''0.126'
$ws.Range('E35').Value = '  -1.71%  '
$ws.Range('E36').Value = '  -2.16%  '
$ws.Range('D37').Value = '''4.96'
$ws.Range('E37').Value = '  +2.71%  '
$ws.Range('D38').Value = '''14.35'
$ws.Range('E38').Value = '  +3.93%  '
$ws.Range('E39').Value = '  -2.54%  '
$ws.Range('B40').Value = 'LidoDAOToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D40').Value = '''2.19'
$ws.Range('E40').Value = '  -7.85%  '
$ws.Range('B41').Value = 'THORChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D41').Value = '''5.87'
$ws.Range('E41').Value = '  -0.72%  '
$ws.Range('D42').Value = '''112.83'
$ws.Range('E42').Value = '  +8.43%  '
$ws.Range('E43').Value = '  -7.54%  '
$ws.Range('D44').Value = '''60.90'
$ws.Range('E44').Value = '  -2.25%  '
$ws.Range('E45').Value = '  -5.84%  '
$ws.Range('D46').Value = '''0.0993'
$ws.Range('E46').Value = '  -4.37%  '
$ws.Range('D47').Value = '''0.997'
$ws.Range('E47').Value = '  -0.43%  '
$ws.Range('D48').Value = '''1.14'
$ws.Range('E48').Value = '  -4.48%  '
$ws.Range('E49').Value = '  -2.00%  '
$ws.Range('D50').Value = '''4.31'
$ws.Range('E50').Value = '  -11.75%  '
$ws.Range('E51').Value = '  +12.44%  '
